# Swap the contents of columns B..AD between row 11 and row 12
# (column A, the running index, stays the same for each row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 11
$row2 = 12

# Columns B (2) through AD (30)
$firstCol = 2
$lastCol = 30

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $cell1 = $ws.Cells.Item($row1, $col)
    $cell2 = $ws.Cells.Item($row2, $col)

    $val1 = $cell1.Value2
    $val2 = $cell2.Value2

    $cell1.Value2 = $val2
    $cell2.Value2 = $val1
}
